$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E12").Value = 22
$ws.Range("E15").Value = 73
$ws.Range("E16").Value = 268
$ws.Range("E18").Value = 75
